$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2, D2, E2 values removed (deleted), C2 value updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.0322458536475985
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: updated values
$ws.Range("B3").Value = 5.3604036204009624
$ws.Range("C3").Value = 5.7505152824675134
$ws.Range("D3").Value = 6.8392903583500395
$ws.Range("E3").Value = 2.5907347027081613

# Selection change to reflect the new active range
$ws.Range("B1:E3").Select()
